$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set number-format to text for the population column so "1.600" etc. stay as literal text
$popRange = $ws.Range("D24:D33")
$popRange.NumberFormat = "@"

$ws.Range("B24").Value = "Королевство Италия"
$ws.Range("C24").Value = "Государство, феодальное, монархическое. Зародилось на терр. совр. Италии в результате отделения от Восточно-Франкского Королевства. Вошло в состав Германского Королевства."
$ws.Range("D24").Value = "1.600"
$ws.Range("E24").Value = "Павия"

$ws.Range("B25").Value = "Королевство Нортумбрия"
$ws.Range("C25").Value = "Одно из семи королевств так называемой англосаксонской гептархии, которое возникло на севере Британии.`nНортумбрия — это англосаксонское государство, образовавшееся в результате объединения в 655 году Берниции и Дейры. После объединения оно просуществовало вплоть до 867 года, когда было захвачено викингами."
$ws.Range("D25").Value = "0.950"
$ws.Range("E25").Value = "Бамборо"

$ws.Range("B26").Value = "Королевство Уэссекс"
$ws.Range("C26").Value = "Англосаксонское королевство на юге Великобритании, одно из семи королевств так называемой англосаксонской гептархии, основанное саксами в начале VI века в ходе англосаксонского завоевания Британии, и ставшее впоследствии частью Англии."
$ws.Range("D26").Value = "0.870"
$ws.Range("E26").Value = "Винчестер"

$ws.Range("B27").Value = "Королевство Мерсия"
$ws.Range("C27").Value = "Одно из семи королевств так называемой англосаксонской гептархии. `nРасполагалось в долине реки Трент на западе центральной Англии "
$ws.Range("D27").Value = "0.970"
$ws.Range("E27").Value = "Тамуэрт"

$ws.Range("B28").Value = "Данелаг"
$ws.Range("C28").Value = "Территория в северо-восточной части Англии, отличавшаяся особыми правовой и социальной системами, унаследованными от норвежских и датских викингов, завоевавших эти земли в IX веке."
$ws.Range("D28").Value = "1.070"
$ws.Range("E28").Value = "Данло"

$ws.Range("B29").Value = "Королевство Нижняя Бургундия"
$ws.Range("C29").Value = "Одно из государств, образовавшихся в процессе распада империи Карла Великого. Королевство включало юго-восточную часть современной Франции."
$ws.Range("D29").Value = "0.700"
$ws.Range("E29").Value = "Вьенн"

$ws.Range("B30").Value = "Королевство Верхняя Бургундия"
$ws.Range("C30").Value = "Одно из государств, образованных после распада империи Карла Великого. Королевство включало территорию западной части современной Швейцарии, Франш-Конте и Шабле. "
$ws.Range("D30").Value = "0.820"
$ws.Range("E30").Value = "Женева"

$ws.Range("B31").Value = "Пражское Княжество"
$ws.Range("C31").Value = "Западнославянское государство X—XII веков, располагавшееся на месте современной Чехии. Правящая династия — Пржемысловичи."
$ws.Range("D31").Value = "0.315"
$ws.Range("E31").Value = "Прага"

$ws.Range("B32").Value = "Империя Ляо"
$ws.Range("C32").Value = "Государство киданей (монголоязычных кочевников), которое занимало просторы Северо-Восточного Китая с момента основания киданьской государственности племенным вождём Абаоцзи"
$ws.Range("D32").Value = "3.800"
$ws.Range("E32").Value = " Шанцзин"

$ws.Range("B33").Value = "Германское Королевство"
$ws.Range("C33").Value = "Государство, феодальное, монархическое. Зародилось на терр. совр. Германии в результате переименования Восточно-Франкского Королевства."
$ws.Range("D33").Value = "5.000"
$ws.Range("E33").Value = "Регенсбург"

# Restore default (General) style on the population cells now that the text is locked in
$popRange.Style = "Normal"

# Leave the active cell where the author left it after the last edit
$ws.Range("C34").Select()
